$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet from "0422" to "0423"
$ws.Name = "0423"

# Fill in the newly-added column J values (rows 2-20)
$ws.Range("J2").Value  = 230
$ws.Range("J3").Value  = 29
$ws.Range("J4").Value  = 51
$ws.Range("J5").Value  = 35
$ws.Range("J6").Value  = 48
$ws.Range("J7").Value  = 35
$ws.Range("J8").Value  = 0
$ws.Range("J9").Value  = 74
$ws.Range("J10").Value = 213
$ws.Range("J11").Value = 29
$ws.Range("J12").Value = 24
$ws.Range("J13").Value = 8
$ws.Range("J14").Value = 47
$ws.Range("J15").Value = 53
$ws.Range("J16").Value = 14
$ws.Range("J17").Value = 0
$ws.Range("J18").Value = 19621
$ws.Range("J19").Value = 23336
$ws.Range("J20").Value = 7814

# Update the active selection to J21 (where the user clicked next)
$ws.Range("J21").Select()
